$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "i"
$ws.Range("B17").Value = "d"
$ws.Range("B17").Select()
